# Update "想去人数" (number of people wanting to attend) figures that changed
# between crawls, for both the "展览" (Exhibition) sheet and the combined
# "全部类型" (All Types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1039
$ws1.Range("F5").Value = 2826
$ws1.Range("F7").Value = 235
$ws1.Range("F11").Value = 101
$ws1.Range("F12").Value = 2661
$ws1.Range("F13").Value = 861

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1039
$ws4.Range("F6").Value = 2826
$ws4.Range("F8").Value = 235
$ws4.Range("F13").Value = 101
$ws4.Range("F14").Value = 2661
$ws4.Range("F15").Value = 861
